$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-17: A=text, B=x, C=y, D=width, E=height
$rows = @(
    @{r=2;  a="February ";   b=908; c=644;    d=83;  e=23},
    @{r=3;  a="Revolution "; b=152; c=669;    d=96;  e=23},
    @{r=4;  a=". ";          b=248; c=669;    d=11;  e=23},
    @{r=5;  a="The ";        b=259; c=669;    d=39;  e=23},
    @{r=6;  a="second ";     b=298; c=669;    d=68;  e=23},
    @{r=7;  a="was ";        b=366; c=669;    d=40;  e=23},
    @{r=8;  a="the ";        b=406; c=669;    d=32;  e=23},
    @{r=9;  a="October ";    b=438; c=669;    d=74;  e=23},
    @{r=10; a="Revolution."; b=512; c=669;    d=96;  e=23},
    @{r=11; a="czarist ";    b=501; c=710.4;  d=61;  e=23},
    @{r=12; a="government."; b=562; c=710.4;  d=107; e=23},
    @{r=13; a="Russian ";    b=819; c=1324.4; d=75;  e=23},
    @{r=14; a="Civil ";      b=894; c=1324.4; d=42;  e=23},
    @{r=15; a="War. ";       b=152; c=1349.4; d=45;  e=23},
    @{r=16; a="Soviet ";     b=491; c=1349.4; d=59;  e=23},
    @{r=17; a="Union ";      b=550; c=1349.4; d=56;  e=23}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
}
